$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-08-21 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-22 Tuesday", 2)

$d.Content.Find.Execute("41×22=", $true, $false, $false, $false, $false, $true, 1, $false, "75×59=", 2)
$d.Content.Find.Execute("14×50=", $true, $false, $false, $false, $false, $true, 1, $false, "97×32=", 2)
$d.Content.Find.Execute("61×66=", $true, $false, $false, $false, $false, $true, 1, $false, "75×87=", 2)
$d.Content.Find.Execute("44×26=", $true, $false, $false, $false, $false, $true, 1, $false, "15×47=", 2)
$d.Content.Find.Execute("41×25=", $true, $false, $false, $false, $false, $true, 1, $false, "20×33=", 2)

$d.Content.Find.Execute("87×52=", $true, $false, $false, $false, $false, $true, 1, $false, "37×62=", 2)
$d.Content.Find.Execute("15×57=", $true, $false, $false, $false, $false, $true, 1, $false, "64×80=", 2)
$d.Content.Find.Execute("28×13=", $true, $false, $false, $false, $false, $true, 1, $false, "53×99=", 2)
$d.Content.Find.Execute("36×28=", $true, $false, $false, $false, $false, $true, 1, $false, "29×20=", 2)
$d.Content.Find.Execute("94×91=", $true, $false, $false, $false, $false, $true, 1, $false, "95×96=", 2)

$d.Content.Find.Execute("21×89=", $true, $false, $false, $false, $false, $true, 1, $false, "14×49=", 2)
$d.Content.Find.Execute("64×85=", $true, $false, $false, $false, $false, $true, 1, $false, "92×49=", 2)
$d.Content.Find.Execute("16×62=", $true, $false, $false, $false, $false, $true, 1, $false, "88×12=", 2)
$d.Content.Find.Execute("73×53=", $true, $false, $false, $false, $false, $true, 1, $false, "21×51=", 2)
$d.Content.Find.Execute("74×14=", $true, $false, $false, $false, $false, $true, 1, $false, "86×38=", 2)

$d.Content.Find.Execute("94×85=", $true, $false, $false, $false, $false, $true, 1, $false, "97×97=", 2)
$d.Content.Find.Execute("35×19=", $true, $false, $false, $false, $false, $true, 1, $false, "48×96=", 2)
$d.Content.Find.Execute("45×55=", $true, $false, $false, $false, $false, $true, 1, $false, "89×64=", 2)
$d.Content.Find.Execute("19×33=", $true, $false, $false, $false, $false, $true, 1, $false, "18×77=", 2)
$d.Content.Find.Execute("15×55=", $true, $false, $false, $false, $false, $true, 1, $false, "63×89=", 2)

$d.Content.Find.Execute("64×75=", $true, $false, $false, $false, $false, $true, 1, $false, "18×13=", 2)
$d.Content.Find.Execute("78×35=", $true, $false, $false, $false, $false, $true, 1, $false, "14×11=", 2)
$d.Content.Find.Execute("68×84=", $true, $false, $false, $false, $false, $true, 1, $false, "73×32=", 2)
$d.Content.Find.Execute("71×98=", $true, $false, $false, $false, $false, $true, 1, $false, "73×96=", 2)
$d.Content.Find.Execute("53×14=", $true, $false, $false, $false, $false, $true, 1, $false, "61×97=", 2)
